$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(3408,45565,24606.25,35.76,4.21,0.9399999999999999),
    @(3409,45566,24800.82,36.04,4.24,0.93),
    @(3410,45568,24314.98,35.33,4.16,0.95),
    @(3411,45569,24066.17,34.97,4.12,0.96),
    @(3412,45572,23404.7,34.01,4,0.99),
    @(3413,45573,23885.55,34.71,4.09,0.97),
    @(3414,45574,24202.38,35.17,4.14,0.96),
    @(3415,45575,24247.84,35.67,4.26,0.96),
    @(3416,45576,24387.39,35.87,4.28,0.95),
    @(3417,45579,24492.31,36.03,4.3,0.95),
    @(3418,45580,24763.62,36.43,4.35,0.9399999999999999),
    @(3419,45581,24767.27,36.43,4.35,0.9399999999999999),
    @(3420,45582,24460.7,36.01,4.3,0.95),
    @(3421,45583,24475.95,36.03,4.3,0.95),
    @(3422,45586,24116.21,35.37,4.24,0.95),
    @(3423,45587,23171.42,33.62,4.07,0.99),
    @(3424,45588,23460.35,33.94,4.12,0.98),
    @(3425,45589,23412.8,34.49,4.11,0.98),
    @(3426,45590,22898.05,33.66,4.02,1.01),
    @(3427,45593,23173.11,34.35,4.07,0.99),
    @(3428,45594,23348.43,34.6,4.1,0.99),
    @(3429,45595,23594.64,34.97,4.14,0.97),
    @(3430,45596,23867.41,35.13,4.19,0.97),
    @(3431,45597,24114.19,35.49,4.24,0.95),
    @(3432,45600,23639.14,34.76,4.15,0.96),
    @(3433,45601,23740.74,34.92,4.17,0.96),
    @(3434,45602,24260.25,35.72,4.26,0.95),
    @(3435,45603,24078.08,35.43,4.23,0.96),
    @(3436,45604,23674.6,34.76,4.16,0.98),
    @(3437,45607,23391.64,34.46,4.11,0.98),
    @(3438,45608,23091.85,34.02,4.05,0.99),
    @(3439,45609,22408.17,33,3.93,1.02),
    @(3440,45610,22592.62,33.31,3.97,1.02),
    @(3441,45614,22473.49,33.11,3.95,1.02),
    @(3442,45615,22695.24,32.7,3.98,1.03),
    @(3443,45617,22591.57,32.51,3.97,1.03),
    @(3444,45618,22795.6,32.81,4,1.01),
    @(3445,45621,23258.18,33.47,4.08,0.99),
    @(3446,45622,23450.06,33.75,4.12,0.98),
    @(3447,45623,23755.84,34.19,4.1,0.97),
    @(3448,45624,23766.96,34.2,4.08,0.97),
    @(3449,45625,23954.16,34.54,4.12,0.99),
    @(3450,45628,24203.45,34.9,4.17,0.98),
    @(3451,45629,24407,35.2,4.2,0.98),
    @(3452,45630,24625.37,35.51,4.24,0.97),
    @(3453,45631,24830.85,35.81,4.28,0.96),
    @(3454,45632,25034.47,36.1,4.31,0.95),
    @(3455,45635,25081.37,36.17,4.32,0.95),
    @(3456,45636,25151.5,36.27,4.33,0.95),
    @(3457,45637,25246.73,36.41,4.35,0.9399999999999999),
    @(3458,45638,25001.63,36.05,4.3,0.95),
    @(3459,45639,24925.57,35.94,4.29,0.96),
    @(3460,45642,25084.48,36.17,4.32,0.95),
    @(3461,45643,24914.19,35.93,4.29,0.96),
    @(3462,45644,24698.28,35.62,4.25,0.96),
    @(3463,45645,24573.39,35.44,4.23,0.97),
    @(3464,45646,24035.53,34.66,4.14,0.99),
    @(3465,45649,24001.46,34.61,4.13,0.99),
    @(3466,45650,24059.05,34.69,4.14,0.99),
    @(3467,45652,24053.96,34.69,4.14,0.99),
    @(3468,45653,24088.9,34.74,4.15,0.99),
    @(3469,45656,23939.99,34.52,4.12,0.99),
    @(3470,45657,24106.01,34.68,4.16,0.98),
    @(3471,45658,24350.79,35.03,4.21,0.97),
    @(3472,45659,24505.64,35.26,4.23,0.97),
    @(3473,45660,24445.73,35.17,4.22,0.97),
    @(3474,45663,23664.3,34.05,4.09,1),
    @(3475,45664,23983.03,34.5,4.14,0.99),
    @(3476,45665,23587.74,33.94,4.08,1),
    @(3477,45666,23270.13,33.48,4.02,1.02),
    @(3478,45667,22662.9,32.61,3.92,1.04),
    @(3479,45670,21733.74,31.27,3.76,1.09),
    @(3480,45671,22164.88,31.9,3.83,1.07)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
